$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (who_n) updates for rows 2-128
$ws.Range("A2").Value = 'Priok, Betul, HP, Mutohar, IM, Posisi, Usai'
$ws.Range("A3").Value = 'Dita Agusta, Rumah Kucing, Dita, Sosok'
$ws.Range("A4").Value = 'Polisi, Kepala Seksi, Rusia, Dinas Dangin Pura, Dilansir, Google Maps, Tigawasa'
$ws.Range("A5").Value = 'Naro, Situs Dalem Lumaju, Rizky Febian, Panyocokan, Umum, Kusworo, Teddy, Wati Tresnawati, Polisi, Islam, Edwin, Narkotika, Raden Ranngamantri, Kapolres, Teddy Pardiyana, Mbah Dalem Lumaju Agung, AKBP Edwin Affandi, Kapolresta, Kombes Kusworo Wibowo, Pengakuan CR, Usai, Satuan Reserse Narkoba Polresta, Ratu, CR, Clandestine'
$ws.Range("A6").Value = 'Pak Kasatlantas, Satlantas, Satlantas Polres, Kanit Laka Lantas Polres, SIM, Badruzzaman, AKP Badruzzaman'
$ws.Range("A7").Value = 'Turus, Sisadani, Celakanya, Paulo, Jose Paulo Ximenes, Toyota Innova, Kanit Laka Lantas Sat Lantas Polres, Innova, Sunan, Sunan Derajat, Darungan, RS, Ipda Wiki Mulyono, Makam Sunan, Neifa, Ayla, AK Zamzami, Wiji'
$ws.Range("A8").Value = 'Petugas Pelabuhan Gilimanuk, Dany, Misalkan, Pelabuhan Gilimanuk'
$ws.Range("A9").Value = 'Binanga, Iptu Junaid, Faldi, Taufik, Idris, Junaid, Muhlis, Polresta'
$ws.Range("A10").Value = 'Bandara Mopah, KNKT, Otban Wilayah X, Corporate Communications Strategic, Supriyadi, Danang Mandala Prihantoro, Lion Air'
$ws.Range("A11").Value = 'AKP Wito, S, Polisi, RBP, A, Kapolsek, Wito'
$ws.Range("A12").Value = 'Sopir, Toyota Innova, Innova, Polres, Polda, A6, Selvi, Yudi Junadi, Yudi, Sugeng, Audi A8, Kompol D, Audi Hitam, Audi, Nur, Supranatural Wowon, Sugeng Guruh, Pembunuhan Wowon, Bantah Tabrak Selvi'
$ws.Range("A13").Value = 'Yohana Rosario Ayunda, Abdullah, Kanit Laka Sat Lantas Polres, BeAT, Sidotopo Sekolahan, Korban, Ngebut, Warsodo, Honda BeAT, Legundi, Iptu Wiji, Wiji'
$ws.Range("A14").Value = 'Kompol Maulana Jali, AKP Agus Suwito, Kanit Laka Lantas Polres, Petamburan, Kasat Lantas Polres, Letjen S Parman, Ali Subchan, Maulana, Duren, Agus'
$ws.Range("A15").Value = 'Usep Supelita, Camat, Ujunggenteng, Sahid, Kades Ujunggenteng Sahid Siam, Usep'
$ws.Range("A16").Value = 'Wiwi, Tol Jagorawi, MT Haryono, Macet, Cawang'
$ws.Range("A17").Value = 'Aipda Benyamin, Kapolda NTT Irjen Johnis Asadoma, Kabid Humas Polda NTT Kombes Pol Ariasandy, Johnis, RSUD Waikabubak, Propam Polres, Daya, Kabid Humas, Polres, Aipda BBA, RS, Aipda Benyamin Anamesa'
$ws.Range("A18").Value = 'Sopir, Sumarjaya, Hamdan, Kasi Humas Polres, Sontak, Hino, Personel, Nyoman Putrawan, Suwug, TKP, Izuzu, AKP I Gede Sumarjaya, RSUD, Pancasari'
$ws.Range("A19").Value = 'Koster, Gunung, I Wayan Koster'
$ws.Range("A20").Value = 'Jasa, Petugas Jasa, Kesulitannya, Pgs, Aprimon, Cawang'
$ws.Range("A21").Value = 'Humas, Azhar Zaki Assjari, Kurnia Rozabi, Lor, Mina, KA, Masinis KA, Petugas, KA Tawangalun, Kereta Api, Kasuran, Muhammad Said, Zaki, Kiai Wachid Hasyim, Niman, Honda Vario'
$ws.Range("A22").Value = 'Muntasir, Sukadamai, Moti Toi, Tolo Kalo, Korban, RSUD, Mustanadi, Manggelewa'
$ws.Range("A23").Value = 'Jogja, AKP Timbul Sasana Raharjo, CCTV, Parsinah, BMW, Shalimaar Style Textile, Timbul, TKP'
$ws.Range("A24").Value = 'Jenis REV, Nagri Kaler, Pratu AA, Pasar, Pemotor, Koncara, SPBU Usman Kebon Kolot, Rantis, TNI, Endoy, Korban, Nagrikaler, Lukman, FindArt Perkebunan Teh, RSUD Bayu Asih'
$ws.Range("A25").Value = 'HIV, ODHA, BUMN, LGBT, Rini, Susah, ARV, PNS, Pengobatan ARV, LPDP ODHA, NTT, AIDS, Rini Maria Denurmin'
$ws.Range("A26").Value = 'Kompol I Wayan Swastika, Minibus, Solihin, Sontak, Gilimanuk, M Suryono, Swastika, Kapolsek'
$ws.Range("A27").Value = 'Semanggi, Kasubdit Gakkum Ditlantas Polda, Kompol Jhoni Eka Putra, Gatot Subroto, Jhoni, Truk, Arif, Udah, S'
$ws.Range("A28").Value = 'Kendaraan, Pospol Tol, Amirudin, Petugas Patroli Jasa, Kedoya'
$ws.Range("A29").Value = 'Iptu Muhammad Rony, Pasir Mandoge, Sofia Nabila, Rony'
$ws.Range("A30").Value = 'Humas Polres, Sumarjaya, Hamdan, Ihda Niswafus Solihah, Nahas, Sontak, Hino, Dinas Buyan, Korban, Suwug, Izuzu, AKP I Gede Sumarjaya, RSUD, Pancasari'
$ws.Range("A31").Value = 'Tambun, Saudara AWK, Truk, Kompol Argadija Putra, Korban, Abah, Tewas Kasat Lantas Polres, AWK, Arga'
$ws.Range("A32").Value = 'Pempem, Petugas Unit Laka Lantas Polsek, Iptu Pempem, Korban, F'
$ws.Range("A33").Value = 'AKP Imam Imam Sayfudin Rodji, Penarikan, Gebangan, Wriginnanom, Imam, Rizal Pratama, TKP, Tol, RSUD, Asta Berris, Abdul Rasek'
$ws.Range("A34").Value = 'Sangarejo, Kasat Reskrim Polres, Farouk, Ponpes, AKP Farouk Ashadi Haiti, Pihaknya'
$ws.Range("A35").Value = 'Sudirman, Diponegoro, Mobil Mitsubishi Xpander, Kompol Ginanjar Fitriadi, RSUD Dr Pirngadi, Baru, Sumatera, Wali, Kapolsek, SIM A, Bobby Nasution, Xpander'
$ws.Range("A36").Value = 'Arga, Tambun, Kasat Lantas Polres, Kompol Argadija Putra, S, RA'
$ws.Range("A37").Value = 'Tambun, Saudara AWK, Truk Ditinggal Sopir Kecelakaan, MRM, Jatimulya, Truk, Kompol Argadija Putra, Korban, Saudari J, AWK, Arga'
$ws.Range("A38").Value = 'Pingky, Tambun, Lalin, Pintu Tol'
$ws.Range("A39").Value = 'Laut Siti Fatimah, Siti, Wowon, Polisi, AKBP Indrawieny Panjiyoga, TKW, Nonik, Duloh, Noneng'
$ws.Range("A40").Value = 'Jogja, AKP Satrio Bagus Wira Wicaksana, Satrio, Plt Kasat Lantas Polresta, Ram Jack, Kopling, Honda CRV, Sajarod, Comby Cutter, Tim SAR, SPBU Armada, Heru Suhartanto, Kepala Kantor SAR, AKBP Mochammad Sajarod Zakun'
$ws.Range("A41").Value = 'Mako Polres, Klinik Pratama, EFK, RS Tripat, Kasat Lantas Polres, ANA, Iptu Agus Rachman, Dasan Geres, EST, Agus, NTB'
$ws.Range("A42").Value = 'Beras, Kapolsek Cikarang, Kompol Sutriesno, Akbar, Triesno, PT, Korban, Polsek Cikarang, Cikarang'
$ws.Range("A43").Value = 'Pegunungan, Penerbad, Selamat, Juinta, KKB'
$ws.Range("A44").Value = 'Khudori, RSUP Kariadi, Ipda Moh Fahrudin, Rozikin, Panit Lantas Polrestabes'
$ws.Range("A45").Value = 'Kompol Maulana Jali, AS, Satlantas, Aipda Satul Bahri, Kasat Lantas Polres, Jakbar, Maulana, Duren, Aipda Bayu Ahmadi, Ipda'
$ws.Range("A46").Value = 'Suteja, Nahas, Kanit Laka Lantas Sat Lantas Polres, Ipda Wiji Mulyono, Darungan, RS, RSUD Ibnu Sina, RS Ibnu Sina, Makam Sunan, Luka, Turus, Kanit Laka Lantas Ipda Wiji Mulyono, Wiji'
$ws.Range("A47").Value = 'Cakra, Kepala Kantor Basarnas, Gede Darmada, Diamond Beach, Koordinator Unit Siaga SAR, Nileash, Austria, Putu Cakra, Penida, SAR'
$ws.Range("A48").Value = 'Kominfo, Dishub, Kapolres, Ekspresi Lapangan Sempur, Satpol PP, Bismo, Kombes Bismo, APAR, Curhat, Pak'
$ws.Range("A49").Value = 'Pelaku, Kapolres, Polisi, Umi Kalsum, AKBP Wiraga Dimas Tama, Wiraga, Motifnya, Anang Budi'
$ws.Range("A50").Value = 'Pasiripis, Pojok, Andi, Welli, Ujunggenteng, Daniel, RS Hermina, Daniel Muttaqien Syafiudin, U, Salakopi, Cipendeuy, Evan'
$ws.Range("A51").Value = 'RS Siti Khadijah, Muchtar, AR Sikakum, Polisi, Kanit Gakkum Satlantas Polrestabes, Jalinsum, Ilir, Raju, Satlantas Polrestabes, RS, TKP, Sumatera, Soekarno Hatta, SP, Iptu AR Sikakum'
$ws.Range("A52").Value = 'Kasat Lantas Polres Pangkep AKP Ida Ayu Made Ari, Tana Toraja, Pangkep, Mobil, RS, Ida'
$ws.Range("A53").Value = 'Kapolres Pakpak Bharat, Rocky, Pakpak Bharat, AKBP Rocky Marpaung, Sumatera, RSUD'
$ws.Range("A54").Value = 'Niko, Pelaku, Asalabuh, Honda Brio RS, Kapolres, Polisi, Kasat Reskrim Polres, Warna Kuning, Jasa, Tinder, AKP Nikolas Bagas Yudhi Kurnia, Catur, Korban, Resmob Polres, SA, AKBP Catur Cahyono Wibowo, Banyuwang'
$ws.Range("A55").Value = 'Sampe, Sarangan, Indah, Allah SWT, Ente, Sugeng, Islam, Habib Usman, Mentalnya'
$ws.Range("A56").Value = 'Kasi Keselamatan Berlayar, Riad, Jembatan Mahkota, Jembatan Mahakam'
$ws.Range("A57").Value = 'Muhammad Ziyad Wijaya, Melki, Ziyad, Melkianus Kotta, Korban, Kutai, SAR, Melkianus'
$ws.Range("A58").Value = 'Lifting Bag, Rescue Rams, Dedik Irianto, Indonesia, Dive Communicator, Eri Cahyadi, Heavy Duty Rescue, Dinas Pemadam Kebakaran, Underwater Drone, Wali, Kepala DPKP, Rescue Spreaders, Mobil Heavy Duty Rescue, Dedik'
$ws.Range("A59").Value = 'Sekaki, Mobil Xenia, Ferdiansyah Prasetya Husada, Satlantas Polresta, Kompol Birgitta Atvina Wijayanti, Korban, Muhammad Alfikri, TKP, Birgitta, Polresta'
$ws.Range("A60").Value = 'Pantai Kelingking, Ni Made Sulistiawati, Dispar, Pemkab, Kepala Dispar, Sulistiawati, One Gate One Destination, Angel Bilabong, Pemerintah, Pemerintah Daerah, Penida, Kelingking'
$ws.Range("A61").Value = 'UPT Transportasi Mamminasata Dishub, Teman Bus, Andi Nur Diyana, Diyana, Armada Teman Bus Trans Maminasata, Sinar, Armada Teman Bus Trans Mamminasata, Boddia, Diana, Fachrul, Husni Mubarak, Panakukang Square'
$ws.Range("A62").Value = 'SIM, Kabid Humas Polda, Satake, Kombes Stefanus Satake Bayu Setianto, Polres, WNA, Polda, Bidang Humas Polda, Polresta'
$ws.Range("A63").Value = 'Boeing, Indonesia, Ahyudin Didakwa Gelapkan Dana, Jaksa, Ampera, Ahyudin, ACT'
$ws.Range("A64").Value = 'Gentong, Anaga, Tasik, Saputro, Satuan Lalu Lintas Polres, Mobil, Lingkar Gentong, Pengemudi, AKP Anaga Budiharso'
$ws.Range("A65").Value = 'Humas, Azhar Zaki Assjari, Lor, Google News, KA, Masinis KA, Petugas, Kereta Api, Kasuran, Zaki, Sutrisno'
$ws.Range("A66").Value = 'MUI, Khofifah, EWS, Irjen Pol Toni Harmanto, Kapolda, Gubernur, Pak Kapolda, Senada'
$ws.Range("A67").Value = 'Campakamekar, Inafis, Direskrimum Polda, Labfor Mabes Polri, TKP, Dedi, Polri, Kadiv Humas Polri Irjen Dedi Prasetyo, Polda, Puslabfor Polri, Kombes K Yani Sudarto'
$ws.Range("A68").Value = 'Masruri, Kahla Anisa, Bupati, Kamulan, Baznas, Deni Riyani, Deni'
$ws.Range("A69").Value = 'Kanit Gakkum Satlantas Polres, Chomsun, Wardatul, Anwar, Saksi Mata, Petugas, Anang, Korban, RSUD, Iptu Anang Setiyanto, Bulurejo'
$ws.Range("A70").Value = 'Sopir, Semeru Putra Transindo, Kemijen, Sarangan, Lurah, Khoirul, Bus, Rombongan, Lawu Green Forest, Duka, Kanti Lestari, AKP Trifona Situmorang'
$ws.Range("A71").Value = 'Kalurahan Jambidan, Jambidan, Polda'
$ws.Range("A72").Value = 'Kapolres, Polisi, Satlantas Polres, Angkot Sinar Murni, Bettes Manurung, AKBP Ronald FC Sipayung, Sumatera, Panribuan, Naga Naga Nagori'
$ws.Range("A73").Value = 'Ikram Saputra, Nahas, Nursalim, Ikram, Panitia'
$ws.Range("A74").Value = 'Selengkapnya, Posko Operasi Zebra, Operasi Zebra, AKBP Made Suarjana, Suarjana, Ditlantas Polda'
$ws.Range("A75").Value = 'Petugas, Komandan Regu Polsuska Kisaran, Sei Dadap, Tri Rahmad Hidayanto, Tri Rahmad'
$ws.Range("A76").Value = 'CCTV, Siswi SMA, Feriza, RR, Korban, Kasat Lantas Polresta Palangkaraya Kompol Feriza Winanda Lubis, Palangkaraya, RA, Jekan, Usai'
$ws.Range("A77").Value = 'Angga, Rengas, Satlantas Polres Muaro, Jaluko, Petugas, AKP Angga Luvyanto, Muaro, Ilham, Rumah Sakit, Kantor Lantas Muaro, Kasat Lantas Polres Muaro'
$ws.Range("A78").Value = 'ESBW, Gakkum, Polisi, Dirlantas Polda, UI, Joko, Kompol Joko Sutriono, TKP, Kombes Latif Usman, Latif, Pemotor, Polri, SOP'
$ws.Range("A79").Value = 'Direktur Lalu Lintas Kepolisian Daerah, Polda, Budi Wahono, Pajero, Ahad, Ajun Komisaris Besar Polisi, Tubuh Hasya, Hasya, UI, WhatsApp, Fadli Zon, X, Komisaris Besar Latif Usman, Mohammad Hasya Athallah Saputra, Polri, Umum Partai Gerindra, DPR RI'
$ws.Range("A80").Value = 'Salaf, Iptu Merdhania Pravita Shanti, Kasi Humas Polres, Lurah, Pasar Besar, Mulyono, Pondok, Muhammad Said, Rejosolor'
$ws.Range("A81").Value = 'Polisi, Kapolres, Lion Air'
$ws.Range("A82").Value = 'Studio Ghibli, Tatsuo Kusakabe, Lune, Haru Yoshioka, Janice Quatlane, Kiki, Kucing, Janice, Melina, Asitaka, Professor Layton, Chihiro, Momo, Lady Eboshi, Diana Wynne Jones, Ponyo, Sistemnya, Eternal Diva, Bibi, Luke, Pokemon, Calcifer, Tombo, Arrietty, Sophie, Film Professor Layton, Totoro, Pangeran, Spirited Away, Ha, Shō, Princess Mononoke, Turniphead, Howl, Shio, Jiji, Miyazaki, Ashitaka, Tokyo, Sosuke, Haru, Delivery Service'
$ws.Range("A83").Value = 'Penumpang, Kasat Lantas Polres, AKP Abdul Malik, Toraja, Bojo, Malik'
$ws.Range("A84").Value = 'EP, Terdakwa EP, Tadi, Jaksa Penuntut Umum, Polisi, Angkutan, Dyah, Jaksa Penuntut Umum Dyah Anggraeni, Dede Halim, Dyah Anggraeni'
$ws.Range("A85").Value = 'Ahmad Nasrul, Sarangan, Semeru, Nasrul, Trifonia, AKP Trifonia Situmorang, Gunung Lawu, Kasat Lantas'
$ws.Range("A86").Value = 'Kombes Kombes Zahwani Pandra Arsyad, Kabid Humas Polda, Pandra, Polda'
$ws.Range("A87").Value = 'PT KCIC, Cempakamekar, Corporate Secretary Rahadian Ratry, Rahadian'
$ws.Range("A88").Value = 'AirNav, Listyo Sigit, Helikopter Polairud, Baharkam Polri, Tempo, Bandara Hanandjoeddin, Kepulauan, Bripda Anam, Briptu Lasminto, Kapolri Jenderal Listyo Sigit Prabowo, Kotawaringin'
$ws.Range("A89").Value = 'Umum PaSKI, Direktur Kepesertaan BPJAMSOSTEK, PaSKI, Zainudin, BPJAMSOSTEK, Semoga, Jarwo'
$ws.Range("A90").Value = 'Mitshubishi Light Truck, Sanggrahan, Honda Supra X, Truk, Korban, Kedungombo, Kasatlantas Polres, Koripan, Maryono, Badan, AKP Maryono'
$ws.Range("A91").Value = 'Viktor Santoso Tandiasa, UU Lalu Lintas, Mahkamah Konstitusi, Mahkamah, DPR, MK, Pemohon, Irfan Kamil'
$ws.Range("A92").Value = 'Kanit Gakkum Satlantas Polres, RS PKU Muhammadiyah, Polisi, KN, Ipda Irwan Marviyanto, Irwan, Ismail Joko Sutrisno, TKP, Jetak, Joko, Soehadi Prijonegoro, RSUD, PMI'
$ws.Range("A93").Value = 'Kepuhanyar, Bypass, Kasat Lantas Polres, KA, AKP Bayu Agustyan, KA Gunung Gedagang, Jery Barokah, Damarsi, Bayu, Empu Nala'
$ws.Range("A94").Value = 'Sarangan, Polisi, Trifonia, Kasat Lantas Polres, PO Semeru Putra Transindo, Ridwan, AKP Trifonia Situmorang, Tim Unit Laka Lantas Polres'
$ws.Range("A95").Value = 'Sarangan, Lawu, Semeru, IGD Puskesmas, Dewi, PO Semeru Putra Transindo, Ridwan, Lawu Green Forest, Sayidiman, RSUD'
$ws.Range("A96").Value = 'APBS, Seksi Tertib Syahbandar, Kombes Pol Puji Hendro, Direktur Polairud Polda, Kantor Kesyahbandaran Utama, Madura, Yuliansyah, Perairan, Karang Jamuang, UM, Jembatan Suramadu, Pertamina, BMKG, Puji, SE Pemkot'
$ws.Range("A97").Value = 'Herman, Kanit Gakkum Polres, Ipda Endang Sudrajat, Gunawan, Endang'
$ws.Range("A98").Value = 'Kadek Ariana, Kabid Humas Polda, Manajer Maruti Grup, Sanur, Resa, ABK, Fast, Iwa Express, Isi Kapal, Satake Bayu, Glory, Ariana, Kombes Stefanus Satake Bayu Setianto, Pelabuhan Sanur, Toni, Ketewel, Agus Purnawijaya, KeboIwa Express, Kebocoran Pengelola Boat Pastikan, Wayan Sadra, Nengah Cemeng'
$ws.Range("A99").Value = 'Ketahuilah, Yatsrib, Lantas, Rasulullah, Kitab, Kiamat, Qais, Nabi SAW, Makkah, Tamim, Islam, Asyraath, Laut Yaman, Kristen, Dajjal Saat, Fatimah, Dajjal, Amir, Arab, Nasrani, Thayyibah'
$ws.Range("A100").Value = 'Fandi Achmad Saputra, Gotekan, Cangar, Kasat Lantas Polres, Wadungasih, Eko, Eko Witomo, Mimik Isbandiyah, Rem, AKP M Bayu Agustyan, Bayu, Wadung, Marcelo Yoga, Honda Vario, Mimik'
$ws.Range("A101").Value = 'Lachlan Brian Hunt, Perampokan Kronologi Kecelakaan Sebut, Aussie, Media Australia Lakukan Penggalangan Dana, Mohon, HP, Polisi, Lachlan, Brian Hunt, Kombes Stefanus Satake Bayu Setianto, Hunt, Lachlan Briant Hunt, Satake Bayu, Polda, Australia, RS BIMC Putri'
$ws.Range("A102").Value = 'Jaksa, Ampera, BCIF, Ibnu Khajar, Ahyudin, ACT, Ibnu Khajar Ahyudin'
$ws.Range("A103").Value = 'Sungai, Kornelis, Kepala Basarnas, Petugas Basarnas, Korban, TNI, Humas Basarnas, Polri'
$ws.Range("A104").Value = 'Pulau Kelor, Pulau Kambing, Taka, Wisata KLM Tiana Liveboat, Nadia, HP, Koordinator Pos SAR, Edy, RS Siloam, Ayu Anjani, CV WAM, Cyntia, Edy Suryono, Kuning, Hasan Sadili, Khouw Cyntia Josephine Kosasih, KTP'
$ws.Range("A105").Value = 'Jenazah Meti, Polisi, Satuan Lalu Lintas Kepolisian Resor, Ariasandy, Meti Tulle, Bidang Hubungan Masyarakat Kepolisian Daerah NTT Komisaris Besar Polisi Ariasandy, Meti'
$ws.Range("A106").Value = 'Halte Cawang UKI, Bernadetta, DKI'
$ws.Range("A107").Value = 'AKBP Bismo Teguh Prakoso, Kamtibmas, Bhabinkamtibmas, Bismo, Kapolresta, Rojali, Polri, Polresta'
$ws.Range("A108").Value = 'Siti, DPRD, Wowon, Dadan, Polisi, Dirkrimum Polda Mentro, Siti Fatimah, Dadan Wandiansyah, Nonik, Wowon Cs, Kombes Hengki Haryadi'
$ws.Range("A109").Value = 'Amerika Serikat, Bamsoet, Kolonel Mar Samson Sitohang, DPR RI Puan Maharani, KSAU Marsekal TNI Fadjar Prasetyo, Kepulauan Seribu, Inggris, Operasi TNI, TNI Angkatan Laut, Konstitusi, Kapolri Jenderal Pol Listyo Sigit Prabowo, Warga Kehormatan Satuan Kapal Selam, Kanada, Brevet Kehormatan Anti Teror Aspek Laut, Kehormatan Korps Marinir, KSAL Laksamana TNI Muhammad Ali, World Directory, Australia, Modern Military Warships, Intai Para Amphibi Korps Marinir, Indonesia, Polri, Pulau, Korps Marinir, Umum Partai Golkar, KSAD Jenderal TNI Dudung Abdurachman, Rusia, Marinir, Umum FKPPI'
$ws.Range("A110").Value = 'Brigjen Ramadhan, Bripda Khoirul Anam, Imade Oka Astawa, Kepala Basarnas Babel, Heli P1103, Divisi Humas Polri Brigjen Ahmad Ramadhan, Perairan, Helikopter, Bandara Pondok Cabe, Briptu Moch Lasminto, Kepulauan, Jenazah Bripda Khoirul Anam, Helikopter Polri Jatuh, Polri, Kapolri Jenderal Listyo Sigit Prabowo, Bun'
$ws.Range("A111").Value = 'Bundaran Songgong, Mandalika, Dimas, Korban, Mandalika AKP I Made Dimas Widiantara, Australia, Sirkuit Mandalika, Antony, NTB'
$ws.Range("A112").Value = 'CCTV, Polsek, Kapolres, Leo Dedy, Polres, Samapta, AKBP Leo Dedy Defretes'
$ws.Range("A113").Value = 'Kapolres, Jimmy, RSR, Polisi, SM, Kasat Narkoba Polres, Rio, AKP Jimmy Ridwan Sihite, AKBP Rio Wahyu Anggoro'
$ws.Range("A114").Value = 'Dadang, Iptu Dadang JB, Tasik, Cintabodas, Kanit Lantas Polsek, Toyota Rush'
$ws.Range("A115").Value = 'Laka, Kompol Ropiyani, Truk, Kasatlantas Polresta'
$ws.Range("A116").Value = 'Heru Budi Hartono, BPNT, BBM, Indonesia, DKI, PKH, APBN, Keadilan, Heru Budi, Gubernur DKI, Menteri Sosial'
$ws.Range("A117").Value = 'Polri Irjen Firman Shantyabudi, Korlantas Polri, Nataru, Komisi V DPR, Senayan, Firman, Kepala Basarnas Henri Alfiandi, Menhub Budi Karya Sumadi'
$ws.Range("A118").Value = 'Sat Lantas, Warsidi, Slamet, Mobil, Tulangrejo, Mobilio'
$ws.Range("A119").Value = 'Edy, Koordinator Pos SAR, Syabandar, Polair, Korban, Edy Suryono, NTT, Pulau Papagarang'
$ws.Range("A120").Value = 'RS Santosa, Yani, Korban, Yani Sudarto, TKP'
$ws.Range("A121").Value = 'AKP Martinus Pararuk, Benny Akbar, Tana Toraja, Orang, Martinus, Lemo, Kapolsek, Lakipadada'
$ws.Range("A122").Value = 'Lampu, SNI, Klakson'
$ws.Range("A123").Value = 'ESBW, SP2HP, Hasya, SP3, Marah, Ira'
$ws.Range("A124").Value = 'Kata Fazzli, Fazzli, Pulau Papatheo, KM, Kepulauan Seribu, Pulau Untung, Kepala Kantor SAR'
$ws.Range("A125").Value = 'Korong, Sopir, Puskesmas Pasar Usang, Anai Iptu Manahan Afrianto Simatupang, Toyota Agya, Sumatera, Kapolsek, Nagari Sungai Buluah, Anai, Agya, Basung'
$ws.Range("A126").Value = 'DS, AR, BAP'
$ws.Range("A127").Value = 'Kombes Pol Shinto Silitonga, ETLE, Polda, Direktorat Lalu Lintas Polda, Bidang Humas Polda, Shinto, Ditlantas Polda'
$ws.Range("A128").Value = 'Betul, Dinas Bina, DKI, Mutohar, DKI Syafrin Liputo, Syafrin'

# new_kec/new_kab/new_prov (M/N/O) corrections
$ws.Range("M23").Value = '-'
$ws.Range("N23").Value = '-'
$ws.Range("O23").Value = '-'
$ws.Range("M38").Value = '-'
$ws.Range("N38").Value = '-'
$ws.Range("O38").Value = '-'
$ws.Range("N91").Value = '-'
$ws.Range("O91").Value = '-'
$ws.Range("M94").Value = 'barat'
$ws.Range("N94").Value = 'magetan'
$ws.Range("O94").Value = 'jawa timur'
$ws.Range("O106").Value = '-'
$ws.Range("M112").Value = 'kuta utara'
$ws.Range("N120").Value = '-'
$ws.Range("O120").Value = '-'
$ws.Range("M121").Value = 'makale utara'
$ws.Range("N123").Value = 'banyuasin'
$ws.Range("O123").Value = 'sumatera selatan'
$ws.Range("M125").Value = 'batang anai'
$ws.Range("N125").Value = 'padang pariaman'
$ws.Range("O125").Value = 'sumatera barat'
